# Updates the cryptocurrency price/volume table with the latest scrape.
# Numeric-looking text values are written with a leading apostrophe so the
# engine stores them as text (matching the source data's inline-string /
# Text cell type) instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.641.19"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "2.323.28"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'312.63"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("D6").Value = "'103.16"
$ws.Range("E6").Value = "  +6.44%  "
$ws.Range("D7").Value = "'0.538"
$ws.Range("E7").Value = "  +2.25%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +7.86%  "
$ws.Range("D10").Value = "'36.23"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").Value = "'0.0820"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("D14").Value = "2.679.24"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").Value = "'15.10"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").Value = "2.334.57"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("D17").Value = "'0.816"
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "43.505.05"
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Value = "0.0₃0943"
$ws.Range("E20").Value = "  +4.17%  "
$ws.Range("D21").Value = "'6.18"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").Value = "'68.52"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "'243.17"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("E24").Value = "  +5.59%  "
$ws.Range("D25").Value = "'2.63"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").Value = "'24.80"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("D29").Value = "'37.79"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.19"
$ws.Range("E30").Value = "  +3.34%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'9.70"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("D32").Value = "'167.30"
$ws.Range("E32").Value = "  +3.23%  "
$ws.Range("D33").Value = "'5.36"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'2.53"
$ws.Range("E35").Value = "  +7.03%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'3.12"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.0751"
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("D38").Value = "'17.83"
$ws.Range("E38").Value = "  +3.80%  "
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").Value = "'1.88"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").Value = "'4.35"
$ws.Range("E42").Value = "  +7.82%  "
$ws.Range("D43").Value = "'19.77"
$ws.Range("E43").Value = "  +4.83%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("D46").Value = "1.984.14"
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").Value = "'3.03"
$ws.Range("E47").Value = "  +5.14%  "
$ws.Range("D48").Value = "'9.91"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").Value = "'56.06"
$ws.Range("E49").Value = "  +5.01%  "
$ws.Range("E50").Value = "  +3.04%  "
$ws.Range("E51").Value = "  +7.47%  "

